$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: product name, unit price, and VAT amount
$ws.Range("A2").Value = "YEDEK PARÇA (24V SELENOİD)"
$ws.Range("D2").Value = 650
$ws.Range("H2").Value = 350000
$ws.Range("I2").ClearContents()

# Clear out row 3 entirely (previously held a second line item)
$ws.Range("A3:I3").ClearContents()

# Keep rows 3, 4 and 5 present as empty rows in the saved sheet
$ws.Rows(3).Hidden = $true
$ws.Rows(3).Hidden = $false
$ws.Rows(4).Hidden = $true
$ws.Rows(4).Hidden = $false
$ws.Rows(5).Hidden = $true
$ws.Rows(5).Hidden = $false
